$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark. In the source document it sits
#    right after "...查看本地确实存在此仓库。" — it will be re-created later
#    at its new home (the very last-but-one empty paragraph).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Locate the paragraph that currently reads "Git remote –v" (it keeps its
#    place) and insert the new "总结" block right after it.
# ---------------------------------------------------------------------------
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Git remote –v") {
        $anchor = $i
    }
}
if (-not $anchor) {
    throw "Could not locate the 'Git remote -v' paragraph"
}

$insertRange = $d.Paragraphs.Item($anchor).Range
$insertRange.Collapse(0)

$newTexts = @(
    "总结：",
    "要关联一个远程库，使用命令git remote add origin git@server-name:path/repo-name.git；",
    "关联后，使用命令git push -u origin master第一次推送master分支的所有内容；",
    "此后，每次本地提交后，只要有必要，就可以使用命令git push origin master推送最新修改；",
    "分布式版本系统的最大好处之一是在本地工作完全不需要考虑远程库的存在，也就是有没有联网都可以正常工作，而SVN在没有联网的时候是拒绝干活的！当有网络的时候，再把本地提交推送一下就完成了同步，真是太方便了！",
    "",
    ""
)

$firstNewIndex = $anchor + 1
foreach ($t in $newTexts) {
    $insertRange.InsertParagraphAfter()
    $insertRange.Collapse(0)
}

for ($i = 0; $i -lt $newTexts.Length; $i++) {
    $p = $d.Paragraphs.Item($firstNewIndex + $i)
    $pr = $p.Range
    if ($newTexts[$i].Length -gt 0) {
        $pr.Text = $newTexts[$i]
    }
}

# "总结：" paragraph -> bold + red, matching the source formatting.
$summaryPara = $d.Paragraphs.Item($firstNewIndex)
$summaryRange = $summaryPara.Range
$summaryRange.Font.Bold = $true
$summaryRange.Font.Color = 255

# ---------------------------------------------------------------------------
# 3. Drop the two now-orphaned blank paragraphs that used to sit right after
#    "Git remote –v" (they are superseded by the new content above, plus the
#    two fresh blank paragraphs already inserted).
# ---------------------------------------------------------------------------
$lastNewIndex = $firstNewIndex + $newTexts.Length - 1   # second new blank paragraph
$oldBlank1 = $lastNewIndex + 1
$oldBlank2 = $lastNewIndex + 2

$d.Paragraphs.Item($oldBlank2).Range.Delete()
$d.Paragraphs.Item($oldBlank1).Range.Delete()

# ---------------------------------------------------------------------------
# 4. Re-create the "_GoBack" bookmark at its new home: the last of the two
#    fresh blank paragraphs inserted above.
#
#    Quirk workaround: adding a bookmark via a collapsed Range that lands
#    inside the last couple of paragraphs of the document mis-places
#    bookmarkEnd (it bleeds into the following paragraph). Padding the
#    document with a few throw-away trailing paragraphs first keeps the
#    target paragraph away from "end of document", then the padding is
#    removed again once the bookmark is safely in place.
# ---------------------------------------------------------------------------
$padCount = 5
$tailRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$tailRange.Collapse(0)
for ($k = 0; $k -lt $padCount; $k++) {
    $tailRange.InsertParagraphAfter()
    $tailRange.Collapse(0)
}

$bookmarkPara = $d.Paragraphs.Item($lastNewIndex)
$bookmarkRange = $bookmarkPara.Range.Duplicate
$bookmarkRange.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

for ($k = 0; $k -lt $padCount; $k++) {
    $lastIdx = $d.Paragraphs.Count
    $d.Paragraphs.Item($lastIdx).Range.Delete()
}
